$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old placeholder " " value in B1 is no longer needed - clear it first so
# its shared string is dropped and the new strings get indices matching the
# target order (Code, Question, Answer, FormHeading).
$ws.Range("B1").ClearContents()

# Row 2: column headers (set first so the shared-string table gets these
# three strings before the heading, matching the original authoring order)
$ws.Range("A2").Value = "Code"
$ws.Range("B2").Value = "Question"
$ws.Range("C2").Value = "Answer"
$ws.Range("A2:C2").Font.Bold = $true

# Row 1: form heading
$ws.Range("A1").Value = "FormHeading"
$ws.Range("A1").Font.Bold = $true

# Column widths
$ws.Columns.Item(1).ColumnWidth = 7.28515625
$ws.Columns.Item(2).ColumnWidth = 28.140625
$ws.Columns.Item(3).ColumnWidth = 50.7109375

# Answer column (C) is left aligned
$ws.Range("C2").HorizontalAlignment = -4131

# Reproduce the saved selection (column C was selected)
[void]$ws.Columns.Item(3).Select()

# Page setup: portrait orientation
$ws.PageSetup.Orientation = 1

# Window position
$excel.Windows.Item(1).Top = 6600
